$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is added for this product at row 69. The previous
# row 69 record gets pushed down to row 70 (dimension grows from R69 to R70).
$ws.Rows.Item(69).Insert()

$ws.Range("A69").Value = 11
$ws.Range("B69").Value = "Vega Monumental Concepción"
$ws.Range("C69").Value = "Bíobío"
$ws.Range("D69").Value = 44595
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Chilena(o)"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 100
$ws.Range("K69").Value = 22000
$ws.Range("L69").Value = 23000
$ws.Range("M69").Value = 22500
$ws.Range("N69").Value = "`$/caja 25 kilos"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 900
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
